$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 with the new tender entry ---
$ws.Range("A2").Value = "Supply of Instructor for Basketball Girls CCA Training Programme"
$ws.Range("B2").Value = "19 May 2021`n01:00PM"
$ws.Range("C2").Value = "Ministry of Education - Schools"
$ws.Range("D2").Value = "Services ⇒ Data Entry, Supply of Manpower Services"
$ws.Range("B2").WrapText = $true

# --- Add new row 3 with another tender entry (includes link-style doc entry) ---
$ws.Range("A3").Value = "Provision of Coaching Services for Basketball (Girls) CCA"
$ws.Range("B3").Value = "12 May 2021`n01:00PM"
$ws.Range("C3").Value = "Ministry of Education - Schools"
$ws.Range("D3").Value = "Services ⇒ Data Entry, Supply of Manpower Services"
$ws.Range("B3").WrapText = $true

# Match row height of row 2 (auto-fit wrapped-text rows are 44.25pt tall)
$ws.Rows.Item(3).RowHeight = 44.25
